# Commit: bb fuel upstream can now multipurpose it be used for both biofuels
# and fossil fuels, if either is connected as a fuel type in a factory
# connections. Using lookup names will propagate from the origin unit process.

$wb = $excel.ActiveWorkbook

# --- "bb electricity" sheet: add two new scenario rows ---
$wsElec = $wb.Worksheets.Item("bb electricity")
$wsElec.Range("A9").Value = "birat-tgr-63vpsa-100bio"
$wsElec.Range("B9").Value = 0.32
$wsElec.Range("C9").Value = "charcoal"
$wsElec.Range("A10").Value = "birat-tgr-100vpsa-100bio"
$wsElec.Range("B10").Value = 0.32
$wsElec.Range("C10").Value = "charcoal"

# --- "bb heat" sheet: add two new scenario rows ---
$wsHeat = $wb.Worksheets.Item("bb heat")
$wsHeat.Range("A9").Value = "birat-tgr-63vpsa-100bio"
$wsHeat.Range("B9").Value = 0.8
$wsHeat.Range("C9").Value = "charcoal"
$wsHeat.Range("A10").Value = "birat-tgr-100vpsa-100bio"
$wsHeat.Range("B10").Value = 0.8
$wsHeat.Range("C10").Value = "charcoal"

# --- "bb fuel upstream" sheet: new "biofuel type" and "biomass type" columns ---
$wsFuel = $wb.Worksheets.Item("bb fuel upstream")

# Header row: insert "biofuel type" before the existing "secondary fuel type"
# column, and append a new "biomass type" column. New header cells pick up
# the same bold 12pt Calibri formatting used by the existing header cells.
$wsFuel.Range("C1").Value = "biofuel type"
$wsFuel.Range("D1").Value = "secondary fuel type"
$wsFuel.Range("D1:E1").Font.Bold = $true
$wsFuel.Range("D1:E1").Font.Size = 12
$wsFuel.Range("D1:E1").Font.Name = "Calibri"

# Existing scenario rows gain biofuel type / secondary fuel type / biomass type values
$wsFuel.Range("B5").Value = "coal"
$wsFuel.Range("C5").Value = "charcoal"
$wsFuel.Range("D5").Value = "natural gas"
$wsFuel.Range("E5").Value = "wood"

$wsFuel.Range("D6").Value = "natural gas"
$wsFuel.Range("C6").Value = "charcoal"
$wsFuel.Range("E6").Value = "wood"

$wsFuel.Range("D7").Value = "natural gas"
$wsFuel.Range("C7").Value = "charcoal"
$wsFuel.Range("E7").Value = "wood"

# "biomass type" header is written after the "wood" values above so the
# workbook's shared-string table keeps the same new-string ordering as the
# target file (wood, then biomass type).
$wsFuel.Range("E1").Value = "biomass type"

# New scenario rows for the 100% biofuel variants
$wsFuel.Range("A8").Value = "birat-tgr-63vpsa-100bio"
$wsFuel.Range("B8").Value = "coal"
$wsFuel.Range("C8").Value = "charcoal"
$wsFuel.Range("D8").Value = "charcoal"
$wsFuel.Range("E8").Value = "wood"

$wsFuel.Range("A9").Value = "birat-tgr-100vpsa-100bio"
$wsFuel.Range("B9").Value = "coal"
$wsFuel.Range("C9").Value = "charcoal"
$wsFuel.Range("D9").Value = "charcoal"
$wsFuel.Range("E9").Value = "wood"

# --- "bb biofuel upstream" sheet: add two new scenario rows ---
$wsBiofuel = $wb.Worksheets.Item("bb biofuel upstream")
$wsBiofuel.Range("A8").Value = "birat-tgr-63vpsa-100bio"
$wsBiofuel.Range("B8").Value = "charcoal"
$wsBiofuel.Range("A9").Value = "birat-tgr-100vpsa-100bio"
$wsBiofuel.Range("B9").Value = "charcoal"
